# 24.12.2021 add fun Revision
# Update the ParentId (column A) values for rows 2-6 on the BOM sheet
# from "Test" to the new revision id "B12UA91111215".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

$ws.Range("A2:A6").Value = "B12UA91111215"
